$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 3)
$full = $cell.Range
$start = $full.Start
$r1 = $d.Range($start, $start + 1)
$r1.Cut() | Out-Null
$cell2 = $t.Cell(4, 3)
Write-Host "after cut:" ($cell2.Range.Text -replace "`r","|" -replace [char]7,"^")
$full2 = $cell2.Range
$r2 = $d.Range($full2.Start, $full2.Start)
$r2.Paste() | Out-Null
$cell3 = $t.Cell(4, 3)
Write-Host "after paste:" ($cell3.Range.Text -replace "`r","|" -replace [char]7,"^")
